# The deck's design was changed from the "Integral" theme to the built-in
# "Office Theme" (Design tab -> Themes gallery -> Office Theme). PowerPoint
# rewrites the theme part backing the slide master (ppt/theme/theme1.xml)
# with the new theme's colour scheme (the font scheme and format scheme of
# "Integral" already matched stock Office values, so only the 12 theme
# colours actually change).
#
# The 12-slot DrawingML colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) is reachable from the object model via ThemeColorScheme, whose
# RGB values are plain COM RGB() integers (0xBBGGRR, i.e. R | G<<8 | B<<16).

$p = $ppt.ActivePresentation

$theme = $p.SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

# Target palette = stock PowerPoint "Office Theme" colour scheme.
$officeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501   # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456   # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477   # folHlink 954F72
}

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i]
}

Write-Host "Applied Office Theme colour scheme to the slide master's theme."
